# Update view-count figures (column F) across the workbook's sheets.
# Mirrors the author's "output generated at 456a3b4" refresh of scraped data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 995
$ws1.Range("F11").Value = 1367
$ws1.Range("F13").Value = 472
$ws1.Range("F14").Value = 1656
$ws1.Range("F18").Value = 1404
$ws1.Range("F19").Value = 270
$ws1.Range("F21").Value = 1135
$ws1.Range("F23").Value = 406
$ws1.Range("F24").Value = 16
$ws1.Range("F25").Value = 3529
$ws1.Range("F26").Value = 699
$ws1.Range("F28").Value = 1564

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F11").Value = 27

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 16

# --- Sheet "全部类型" (All Types, union of the above) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 16
$ws4.Range("F15").Value = 995
$ws4.Range("F21").Value = 1367
$ws4.Range("F23").Value = 472
$ws4.Range("F24").Value = 1656
$ws4.Range("F28").Value = 1404
$ws4.Range("F29").Value = 270
$ws4.Range("F32").Value = 27
$ws4.Range("F33").Value = 1135
$ws4.Range("F35").Value = 406
$ws4.Range("F36").Value = 16
$ws4.Range("F37").Value = 3529
$ws4.Range("F38").Value = 699
$ws4.Range("F40").Value = 1564
